$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the vocabulary namespace URIs (test3 -> test2) ---
$ws.Range("B1").Value = "http://purl.org/test2/variables/"
$ws.Range("C3").Value = "http://purl.org/test2/variables/"

# --- Update the title / description placeholders ---
$ws.Range("B10").Value = "Test2 vocabulary"
$ws.Range("B11").Value = "Test for showcase purposes"

# --- Remove the duplicate "dct:creator | Minka" row (row 13). ---
# This shifts every row below it up by one, which also accounts for the
# dimension shrinking from A1:T90 to A1:T89 (the now-superfluous trailing
# blank "vars:" row disappears off the bottom).
$ws.Rows(13).Delete()

# --- After the shift, populate the two vocabulary term rows (now 19 & 20)
#     and the row that used to be an empty "vars:" placeholder (now 21). ---

# Row 19: vars:Test
$ws.Range("A19").Value = "vars:Test"
$ws.Range("B19").Value = "Test"
$ws.Range("E19").Value = "Test for technical setup"

# Row 20: vars:Computerscientist
$ws.Range("A20").Value = "vars:Computerscientist"
$ws.Range("B20").Value = "Computerscientist"
$ws.Range("E20").Value = "person that studied computerscience"
$ws.Range("F20").Value = "vars:Computerscience"

# Row 21: vars:Computerscience
$ws.Range("A21").Value = "vars:Computerscience"
$ws.Range("B21").Value = "Computerscience"
$ws.Range("E21").Value = "studies of computers"
